$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = [double]"1.000000102806806"
$ws.Cells.Item(2, 6).Value = [double]"2.949214424937967E-07"
$ws.Cells.Item(2, 7).Value = [double]"3390733.797959907"
$ws.Cells.Item(2, 8).Value = [double]"2.603370599182205E-147"

$ws.Cells.Item(3, 5).Value = [double]"1.23433784887034"
$ws.Cells.Item(3, 6).Value = [double]"0.3166090886962632"
$ws.Cells.Item(3, 7).Value = [double]"3.898617863287222"
$ws.Cells.Item(3, 8).Value = [double]"0.0006425103851562793"

$ws.Cells.Item(4, 5).Value = [double]"1.232533273701805"
$ws.Cells.Item(4, 6).Value = [double]"0.4287928708670447"
$ws.Cells.Item(4, 7).Value = [double]"2.874425759946871"
$ws.Cells.Item(4, 8).Value = [double]"0.008147289954327742"

$ws.Cells.Item(5, 5).Value = [double]"1.371350092038687"
$ws.Cells.Item(5, 6).Value = [double]"0.424547767205771"
$ws.Cells.Item(5, 7).Value = [double]"3.230143220548413"
$ws.Cells.Item(5, 8).Value = [double]"0.003451021144898397"

$ws.Cells.Item(6, 5).Value = [double]"1.038715776342111"
$ws.Cells.Item(6, 6).Value = [double]"0.3312664567389792"
$ws.Cells.Item(6, 7).Value = [double]"3.135589961529262"
$ws.Cells.Item(6, 8).Value = [double]"0.004349323766392233"

$ws.Cells.Item(7, 5).Value = [double]"1.499984864011523"
$ws.Cells.Item(7, 6).Value = [double]"0.3309360513054669"
$ws.Cells.Item(7, 7).Value = [double]"4.532552008445217"
$ws.Cells.Item(7, 8).Value = [double]"0.0001251695443991448"

$ws.Cells.Item(8, 5).Value = [double]"1.432969991410402"
$ws.Cells.Item(8, 6).Value = [double]"0.5636909536328849"
$ws.Cells.Item(8, 7).Value = [double]"2.542119901295512"
$ws.Cells.Item(8, 8).Value = [double]"0.01759203995199525"

$ws.Cells.Item(9, 5).Value = [double]"1.261901771453531"
$ws.Cells.Item(9, 6).Value = [double]"0.2571949016513472"
$ws.Cells.Item(9, 7).Value = [double]"4.906402744966393"
$ws.Cells.Item(9, 8).Value = [double]"4.748950712240083E-05"

$ws.Cells.Item(10, 5).Value = [double]"0.8905381796509941"
$ws.Cells.Item(10, 6).Value = [double]"0.2076663778679619"
$ws.Cells.Item(10, 7).Value = [double]"4.288311804702515"
$ws.Cells.Item(10, 8).Value = [double]"0.0002356127056390428"

$ws.Cells.Item(11, 5).Value = [double]"0.6464830737655988"
$ws.Cells.Item(11, 6).Value = [double]"0.296419247169613"
$ws.Cells.Item(11, 7).Value = [double]"2.180975358174623"
$ws.Cells.Item(11, 8).Value = [double]"0.03880269101825961"
